$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header cells (Wins, Losses, Ties) in AC1:AE1, copying the
# formatting used by the existing header row (bold, bordered, centered).
$ws.Range("AB1").Copy()
$ws.Range("AC1:AE1").PasteSpecial(-4122)

$ws.Range("AC1").Value = "Wins"
$ws.Range("AD1").Value = "Losses"
$ws.Range("AE1").Value = "Ties"

# Fill in the season record (Wins/Losses/Ties) for every data row.
for ($r = 2; $r -le 45; $r++) {
    $ws.Range("AC$r").Value = 94
    $ws.Range("AD$r").Value = 68
    $ws.Range("AE$r").Value = 0
}
